$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 - this pushes the existing rows 18-27
# down to 19-28 (matching the diff's "shift everything down by one row"
# pattern), and the sheet's used range grows to A1:R28.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly price record.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44664
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100114007
$ws.Range("G18").Value = "Jengibre"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 11000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 11600
$ws.Range("N18").Value = "$/caja 13 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 892
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "Hortaliza"
